$wb = $excel.ActiveWorkbook

$wsProfits = $wb.Worksheets.Item("CRtPaL-profits")
$wsLosses  = $wb.Worksheets.Item("CRtPaL-losses")

# Rename existing "hydrogen" entry to "hydrogen combustion turbine" (shared string rename,
# applies to row 24 on both sheets since they reference the same shared string).
$wsProfits.Range("A24").Value = "hydrogen combustion turbine"
$wsLosses.Range("A24").Value = "hydrogen combustion turbine"

# Update capacity response values for rows 16 and 17 on the profits sheet.
$wsProfits.Range("B16").Value = 0
$wsProfits.Range("B17").Value = 0

# Add the new "hydrogen combined cycle" row (row 25) on both sheets.
$wsProfits.Range("A25").Value = "hydrogen combined cycle"
$wsProfits.Range("B25").Value = 2

$wsLosses.Range("A25").Value = "hydrogen combined cycle"
$wsLosses.Range("B25").Value = 1

# Match B25's number format to B24's (integer display format) on the losses sheet.
$wsLosses.Range("B25").NumberFormat = $wsLosses.Range("B24").NumberFormat

# Apply the new font/alignment formatting (black font color, vertically centered) to the two
# "hydrogen..." rows on both sheets. Build the combined format once on a scratch cell and paste
# it to both sheets, so only a single new font/style entry is produced (matching the source
# formatting) instead of one per sheet.
$scratch = $wsProfits.Range("Z100")
$scratch.Font.Color = 0
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$wsProfits.Range("A24:A25").PasteSpecial(-4122)
$wsLosses.Range("A24:A25").PasteSpecial(-4122)
$scratch.Clear()
